# Update "想去人数" (F column) figures for a handful of events on the
# "展览" and "全部类型" sheets, reflecting refreshed scrape output.

$wb = $excel.ActiveWorkbook

# Sheet "展览": rows 2, 5, 6, 9
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 3455
$wsExpo.Range("F5").Value = 6993
$wsExpo.Range("F6").Value = 2510
$wsExpo.Range("F9").Value = 29

# Sheet "全部类型": rows 2, 6, 7, 10 (shifted by one row vs "展览")
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 3455
$wsAll.Range("F6").Value = 6993
$wsAll.Range("F7").Value = 2510
$wsAll.Range("F10").Value = 29
